$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CIS to CCMS import analysis")

$ws.Range("A2").Value = 45651.0
$ws.Range("B2").Value = "Invisibility potion"

$ws.Range("D3").Value = 2.0

$ws.Range("B5").Value = "Energy increase"

$ws.Range("B6").Value = "Extra energy increase"
